$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "o1980, o1988, o1992b, o2008"
$ws.Range("E4").Value = "o1980, o1988, o1992b"
$ws.Range("E5").Value = "o1980, o1988, o2008"
$ws.Range("E6").Value = "o1980, o1988, o1992b, o2008"
$ws.Range("E7").Value = "o1980, o1988, o1992b"
$ws.Range("E8").Value = "o1980, o1988, o1992b"
$ws.Range("E9").Value = "o1980, o1988, o1992b"
$ws.Range("E10").Value = "o1980, o1988, o2008"
$ws.Range("E11").Value = "o1980, o1988, o1992b"
$ws.Range("E12").Value = "o1980, o1988, o1992b"
$ws.Range("E14").Value = "o1980, o1988"
$ws.Range("E15").Value = "o1980, o1992b"
$ws.Range("E17").Value = "o1980"
$ws.Range("E18").Value = "o1980, o1988, o1992b"
$ws.Range("E19").Value = "o1980, o1988, o1992b, o2008"
$ws.Range("E21").Value = "o1980, o1988"
$ws.Range("B22").Value = "o1980"
$ws.Range("E22").Value = "o1980, o1988"
$ws.Range("E23").Value = "o1980, o1988, o1992b"
$ws.Range("E25").Value = "o1980, o1988"
$ws.Range("B26").Value = "o1980"
$ws.Range("E26").Value = "o1980, o1988, o1992b"
$ws.Range("B27").Value = "o2008"
$ws.Range("E27").Value = "o2008"
$ws.Range("B28").Value = "o2008"
$ws.Range("E28").Value = "o2008"
$ws.Range("B29").Value = "o1988"
$ws.Range("E29").Value = "o1988"
$ws.Range("B31").Value = "o2008"
$ws.Range("E31").Value = "o2008"
$ws.Range("B34").Value = "o1980"
$ws.Range("E34").Value = "o1980, o1988, o1992b"
$ws.Range("E35").Value = "o1980"
$ws.Range("E36").Value = "o1980, o1988, o1992b"
$ws.Range("E37").Value = "o1980"
$ws.Range("B38").Value = "o1980"
$ws.Range("E38").Value = "o1980, o1988"
$ws.Range("E39").Value = "o1980"
$ws.Range("E41").Value = "o1980, o1988"
$ws.Range("B43").Value = "o1988"
$ws.Range("E43").Value = "o1988, o1992b"
$ws.Range("B44").Value = "o1980"
$ws.Range("E44").Value = "o1980, o1988, o1992b"
$ws.Range("B47").Value = "o2008"
$ws.Range("E47").Value = "o2008"
$ws.Range("B48").Value = "o2008"
$ws.Range("E48").Value = "o2008"
